$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the results table (rows 2-14) with the final recomputed values.
# Columns: A=Rank, B=Classifier, C=Features, D=Parameters, E=Accuracy, F=Recall,
#          G=Precision, H=AUC, I=True_Positives, J=True_Negatives,
#          K=False_Positives, L=False_Negatives

$ws.Range("A2").Value = 12
$ws.Range("B2").Value = "XGB"
$ws.Range("C2").Value = "Full"
$ws.Range("D2").Value = "Default"
$ws.Range("E2").Value = 0.9514063697055927
$ws.Range("F2").Value = 0.9545953059881553
$ws.Range("G2").Value = 0.9483547613859229
$ws.Range("H2").Value = 0.9933575380926222
$ws.Range("I2").Value = 8704
$ws.Range("J2").Value = 8682
$ws.Range("K2").Value = 474
$ws.Range("L2").Value = 414

$ws.Range("A3").Value = 9
$ws.Range("B3").Value = "StackingCV"
$ws.Range("C3").Value = "Reduced"
$ws.Range("D3").Value = "Best"
$ws.Range("E3").Value = 0.9533216591879172
$ws.Range("F3").Value = 0.9634788330774292
$ws.Range("G3").Value = 0.9441160666308437
$ws.Range("H3").Value = 0.9931383833972925
$ws.Range("I3").Value = 8785
$ws.Range("J3").Value = 8636
$ws.Range("K3").Value = 520
$ws.Range("L3").Value = 333

$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Stacking (SGD)"
$ws.Range("C4").Value = "Reduced"
$ws.Range("D4").Value = "Best"
$ws.Range("E4").Value = 0.9532669366312794
$ws.Range("F4").Value = 0.9639175257731959
$ws.Range("G4").Value = 0.9436332402834443
$ws.Range("H4").Value = 0.9931375449173695
$ws.Range("I4").Value = 8789
$ws.Range("J4").Value = 8631
$ws.Range("K4").Value = 525
$ws.Range("L4").Value = 329

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Random Forest"
$ws.Range("C5").Value = "Full"
$ws.Range("D5").Value = "Default"
$ws.Range("E5").Value = 0.9504760862427493
$ws.Range("F5").Value = 0.9754332090370695
$ws.Range("G5").Value = 0.9288772845953003
$ws.Range("H5").Value = 0.9922658432218865
$ws.Range("I5").Value = 8894
$ws.Range("J5").Value = 8475
$ws.Range("K5").Value = 681
$ws.Range("L5").Value = 224

$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Voting"
$ws.Range("C6").Value = "Reduced"
$ws.Range("D6").Value = "Best"
$ws.Range("E6").Value = 0.9447302177957754
$ws.Range("F6").Value = 0.9599692915112963
$ws.Range("G6").Value = 0.931368376250266
$ws.Range("H6").Value = 0.9907880523031318
$ws.Range("I6").Value = 8753
$ws.Range("J6").Value = 8511
$ws.Range("K6").Value = 645
$ws.Range("L6").Value = 365

$ws.Range("A7").Value = 8
$ws.Range("B7").Value = "Stacking (SVC)"
$ws.Range("C7").Value = "Reduced"
$ws.Range("D7").Value = "Best"
$ws.Range("E7").Value = 0.9541972200941228
$ws.Range("F7").Value = 0.9598596183373547
$ws.Range("G7").Value = 0.9489320177816328
$ws.Range("H7").Value = 0.9770402396576855
$ws.Range("I7").Value = 8752
$ws.Range("J7").Value = 8685
$ws.Range("K7").Value = 471
$ws.Range("L7").Value = 366

$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "SVM (SVC)"
$ws.Range("C8").Value = "Full"
$ws.Range("D8").Value = "Default"
$ws.Range("E8").Value = 0.9111232279171211
$ws.Range("F8").Value = 0.9393939393939394
$ws.Range("G8").Value = 0.8902564102564102
$ws.Range("H8").Value = 0.9513557870700727
$ws.Range("I8").Value = 868
$ws.Range("J8").Value = 803
$ws.Range("K8").Value = 107
$ws.Range("L8").Value = 56

$ws.Range("A9").Value = 10
$ws.Range("B9").Value = "Vecstack"
$ws.Range("C9").Value = "Reduced"
$ws.Range("D9").Value = "Best"
$ws.Range("E9").Value = 0.9473021779577542
$ws.Range("F9").Value = 0.9787234042553191
$ws.Range("G9").Value = 0.9206643969875168
$ws.Range("H9").Value = 0.9473673814636141
$ws.Range("I9").Value = 8924
$ws.Range("J9").Value = 8387
$ws.Range("K9").Value = 769
$ws.Range("L9").Value = 194

$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "Logistic Regression"
$ws.Range("C10").Value = "Full"
$ws.Range("D10").Value = "Default"
$ws.Range("E10").Value = 0.8752872934223487
$ws.Range("F10").Value = 0.8895591138407546
$ws.Range("G10").Value = 0.8644356815517426
$ws.Range("H10").Value = 0.9349048507357206
$ws.Range("I10").Value = 8111
$ws.Range("J10").Value = 7884
$ws.Range("K10").Value = 1272
$ws.Range("L10").Value = 1007

$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "Decision Tree"
$ws.Range("C11").Value = "Full"
$ws.Range("D11").Value = "Default"
$ws.Range("E11").Value = 0.9251942650760644
$ws.Range("F11").Value = 0.9325509980258828
$ws.Range("G11").Value = 0.9187466234467855
$ws.Range("H11").Value = 0.9252095313414691
$ws.Range("I11").Value = 8503
$ws.Range("J11").Value = 8404
$ws.Range("K11").Value = 752
$ws.Range("L11").Value = 615

$ws.Range("A12").Value = 6
$ws.Range("B12").Value = "Stacking (Logistic)"
$ws.Range("C12").Value = "Reduced"
$ws.Range("D12").Value = "Best"
$ws.Range("E12").Value = 0.9201597898653825
$ws.Range("F12").Value = 0.9250932221978504
$ws.Range("G12").Value = 0.9157529041363587
$ws.Range("H12").Value = 0.9203934883265866
$ws.Range("I12").Value = 8435
$ws.Range("J12").Value = 8380
$ws.Range("K12").Value = 776
$ws.Range("L12").Value = 683

$ws.Range("A13").Value = 2
$ws.Range("B13").Value = "Gaussian Naive-Bayes"
$ws.Range("C13").Value = "Full"
$ws.Range("D13").Value = "Default"
$ws.Range("E13").Value = 0.7587282477837365
$ws.Range("F13").Value = 0.6627549901294143
$ws.Range("G13").Value = 0.8191676833401111
$ws.Range("H13").Value = 0.8552142994174434
$ws.Range("I13").Value = 6043
$ws.Range("J13").Value = 7822
$ws.Range("K13").Value = 1334
$ws.Range("L13").Value = 3075

$ws.Range("A14").Value = 0
$ws.Range("B14").Value = "Bernoulli Naive-Bayes"
$ws.Range("C14").Value = "Full"
$ws.Range("D14").Value = "Default"
$ws.Range("E14").Value = 0.6709532669366313
$ws.Range("F14").Value = 0.6650581267821891
$ws.Range("G14").Value = 0.6720602903690569
$ws.Range("H14").Value = 0.7323666294669059
$ws.Range("I14").Value = 6064
$ws.Range("J14").Value = 6197
$ws.Range("K14").Value = 2959
$ws.Range("L14").Value = 3054
